# Auto-generated Word COM-interop script: translate ContosoLearn Market Research doc to Italian
$d = $word.ActiveDocument

function Set-ParaXml($paraIndex, [string]$bodyXml) {
    $p = $d.Paragraphs($paraIndex).Range
    $r = $d.Range($p.Start, $p.End - 1)
    $pkg = '<?xml version="1.0" encoding="utf-8"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $bodyXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
    $d.Paragraphs($paraIndex).Format.ReadingOrder = 0
}

$body1 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>Ricerca di mercato di ContosoLearn</w:t></w:r>
'@
Set-ParaXml 1 $body1

$body2 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>AdatumLearn: AdatumLearn è una delle principali piattaforme di apprendimento basate sull'intelligenza artificiale, che arricchisce l'eLearning automatizzando un'ampia gamma di attività.</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>È nota per le sue funzionalità di creazione di contenuti e per la tecnologia di apprendimento adattivo.</w:t></w:r>
'@
Set-ParaXml 2 $body2

$body3 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>AdventureLearn: AdventureLearn è un'altra piattaforma di apprendimento basata sull'intelligenza artificiale, che offre esperienze di apprendimento personalizzate e raccomandazioni basate sull'analisi dei dati.</w:t></w:r>
'@
Set-ParaXml 3 $body3

$body4 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>AlpineTraining: AlpineTraining è una piattaforma di apprendimento ottimizzata per l'uso su dispositivi mobili, focalizzata sul microlearning.</w:t></w:r>
'@
Set-ParaXml 4 $body4

$body5 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>Bellows OnDemand: Bellows OnDemand offre una soluzione di apprendimento completa, comprendente strumenti per la creazione di contenuti e la collaborazione sociale.</w:t></w:r>
'@
Set-ParaXml 5 $body5

$body6 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>FabrikamLearning: FabrikamLearning fornisce una suite di piattaforme in grado di rispondere a diverse esigenze di apprendimento.</w:t></w:r>
'@
Set-ParaXml 6 $body6

$body7 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>FirstUp Cards: FirstUp Cards è un'app di apprendimento mobile ideale per gestire training su procedure di sicurezza, conformità, conoscenza di nuovi prodotti o qualsiasi altra tipologia di scenario di apprendimento.</w:t></w:r>
'@
Set-ParaXml 7 $body7

$body8 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>Munson'sLearn: Munson'sLearn supporta le aziende nella formazione di dipendenti, partner e clienti.</w:t></w:r>
'@
Set-ParaXml 8 $body8

$body9 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>LibertyLearn: LibertyLearn è un LMS veloce e progettato per progetti mission-critical.</w:t></w:r>
'@
Set-ParaXml 9 $body9

$body10 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>WoodgroveLMS: WoodgroveLMS è un sistema di gestione dell'apprendimento funzionale ed esteticamente curato, concepito per offrire un'esperienza didattica di alto livello.</w:t></w:r>
'@
Set-ParaXml 10 $body10

$body11 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>NorthwindWorlds: NorthwindWorlds è una soluzione di apprendimento potente, intuitiva e affidabile, adatta sia a singoli individui che a realtà aziendali.</w:t></w:r>
'@
Set-ParaXml 11 $body11

$body12 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>ProsewareLearn: ProsewareLearn è un'azienda di formazione online che propone un'ampia selezione di corsi video per sviluppatori software, amministratori IT e professionisti creativi.</w:t></w:r>
'@
Set-ParaXml 12 $body12

$body13 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>RelecloudLearn: RelecloudLearn è una piattaforma di apprendimento online che offre MOOC (corsi aperti su larga scala), specializzazioni e titoli di studio in numerose discipline.</w:t></w:r>
'@
Set-ParaXml 13 $body13

$body14 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>TreyAcademy: TreyAcademy è una piattaforma di apprendimento online, lanciata nel maggio 2010, rivolta a professionisti e studenti.</w:t></w:r>
'@
Set-ParaXml 14 $body14

$body15 = @'
<w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>Queste piattaforme vantano una notevole presenza sul mercato e sono ampiamente riconosciute per le loro funzionalità basate sull'AI, tra cui esperienze di apprendimento personalizzate, raccomandazioni guidate dai dati e automazione dei processi.</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>Stanno trasformando il panorama di eLearning, sfruttando l'intelligenza artificiale per offrire esperienze di apprendimento più coinvolgenti, gratificanti e su misura per ogni utente.</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="it-IT" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
'@
Set-ParaXml 15 $body15

Write-Output "OK"
